# Update "演出" (performances) sheet: remove the cancelled/removed event row
# (2024-08-23 "北京·《山丘》音乐教父..."), which shifts all following rows
# up by one, and bump the "想去人数" (want-to-go count) on the last
# remaining row.
$wb = $excel.ActiveWorkbook
$wsPerf = $wb.Worksheets.Item(2)

$wsPerf.Rows(4).Delete()

# Column A holds a plain sequential index (0-based, independent of the
# underlying event) -- restore it after the row shift.
for ($r = 4; $r -le 15; $r++) {
    $wsPerf.Cells.Item($r, 1).Value = $r - 1
}

# Last row's want-to-go count increased from 217 to 218.
$wsPerf.Cells.Item(15, 6).Value = 218

# Update "想去人数" (want-to-go count, column F) values across the sheets
# to the latest scraped totals.

# "展览" (exhibitions)
$wsExpo = $wb.Worksheets.Item(1)
$wsExpo.Cells.Item(3, 6).Value = 428
$wsExpo.Cells.Item(4, 6).Value = 169
$wsExpo.Cells.Item(5, 6).Value = 3947
$wsExpo.Cells.Item(7, 6).Value = 2602
$wsExpo.Cells.Item(8, 6).Value = 85
$wsExpo.Cells.Item(9, 6).Value = 3219
$wsExpo.Cells.Item(11, 6).Value = 2345
$wsExpo.Cells.Item(13, 6).Value = 119
$wsExpo.Cells.Item(14, 6).Value = 330
$wsExpo.Cells.Item(15, 6).Value = 471
$wsExpo.Cells.Item(16, 6).Value = 22
$wsExpo.Cells.Item(18, 6).Value = 229
$wsExpo.Cells.Item(20, 6).Value = 315
$wsExpo.Cells.Item(21, 6).Value = 461
$wsExpo.Cells.Item(22, 6).Value = 683
$wsExpo.Cells.Item(23, 6).Value = 1435
$wsExpo.Cells.Item(24, 6).Value = 250
$wsExpo.Cells.Item(26, 6).Value = 1312
$wsExpo.Cells.Item(27, 6).Value = 146
$wsExpo.Cells.Item(28, 6).Value = 170
$wsExpo.Cells.Item(29, 6).Value = 8
$wsExpo.Cells.Item(30, 6).Value = 77
$wsExpo.Cells.Item(31, 6).Value = 4453
$wsExpo.Cells.Item(32, 6).Value = 4352
$wsExpo.Cells.Item(33, 6).Value = 93
$wsExpo.Cells.Item(34, 6).Value = 300
$wsExpo.Cells.Item(37, 6).Value = 1168
$wsExpo.Cells.Item(38, 6).Value = 164
$wsExpo.Cells.Item(40, 6).Value = 505
$wsExpo.Cells.Item(43, 6).Value = 187
$wsExpo.Cells.Item(45, 6).Value = 118
$wsExpo.Cells.Item(46, 6).Value = 46
$wsExpo.Cells.Item(47, 6).Value = 71
$wsExpo.Cells.Item(49, 6).Value = 5

# "本地生活" (local life)
$wsLocal = $wb.Worksheets.Item(3)
$wsLocal.Cells.Item(3, 6).Value = 159
$wsLocal.Cells.Item(4, 6).Value = 2341
$wsLocal.Cells.Item(5, 6).Value = 42

# "全部类型" (all types)
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Cells.Item(3, 6).Value = 159
$wsAll.Cells.Item(4, 6).Value = 428
$wsAll.Cells.Item(5, 6).Value = 169
$wsAll.Cells.Item(6, 6).Value = 3947
$wsAll.Cells.Item(7, 6).Value = 2602
$wsAll.Cells.Item(8, 6).Value = 85
$wsAll.Cells.Item(9, 6).Value = 3219
$wsAll.Cells.Item(12, 6).Value = 2345
$wsAll.Cells.Item(14, 6).Value = 119
$wsAll.Cells.Item(15, 6).Value = 330
$wsAll.Cells.Item(16, 6).Value = 471
$wsAll.Cells.Item(17, 6).Value = 22
$wsAll.Cells.Item(18, 6).Value = 229
$wsAll.Cells.Item(21, 6).Value = 683
$wsAll.Cells.Item(22, 6).Value = 1435
$wsAll.Cells.Item(23, 6).Value = 1312
$wsAll.Cells.Item(24, 6).Value = 146
$wsAll.Cells.Item(26, 6).Value = 77
$wsAll.Cells.Item(29, 6).Value = 4453
$wsAll.Cells.Item(30, 6).Value = 4352
$wsAll.Cells.Item(31, 6).Value = 93
$wsAll.Cells.Item(33, 6).Value = 1168
$wsAll.Cells.Item(34, 6).Value = 164
$wsAll.Cells.Item(38, 6).Value = 505
$wsAll.Cells.Item(44, 6).Value = 187
$wsAll.Cells.Item(45, 6).Value = 118
$wsAll.Cells.Item(46, 6).Value = 46
$wsAll.Cells.Item(47, 6).Value = 71
$wsAll.Cells.Item(49, 6).Value = 218
